$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 12.97098566666667
$ws.Cells.Item(2, 8).Value = 38.91295700000001
$ws.Cells.Item(2, 9).Value = 0.7291028508134716
$ws.Cells.Item(2, 10).Value = 0.7291028508134717
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 47.57896333333333
$ws.Cells.Item(2, 14).Value = 142.73689
$ws.Cells.Item(2, 15).Value = 0.450188452948237
$ws.Cells.Item(2, 16).Value = 0.4501884529482371
$ws.Cells.Item(2, 17).Value = 617.1460514315256
$ws.Cells.Item(2, 18).Value = 5554.314462883731
$ws.Cells.Item(2, 19).Value = 0.328233684447866
$ws.Cells.Item(2, 20).Value = 0.3282336844478661

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 12.97098566666667
$ws.Cells.Item(3, 8).Value = 38.91295700000001
$ws.Cells.Item(3, 9).Value = 0.7291028508134716
$ws.Cells.Item(3, 10).Value = 0.7291028508134717
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 9.829723666666666
$ws.Cells.Item(3, 14).Value = 29.489171
$ws.Cells.Item(3, 15).Value = 0.09300808131111737
$ws.Cells.Item(3, 16).Value = 0.09300808131111739
$ws.Cells.Item(3, 17).Value = 127.5012047876275
$ws.Cells.Item(3, 18).Value = 1147.510843088647
$ws.Cells.Item(3, 19).Value = 0.06781245723262684
$ws.Cells.Item(3, 20).Value = 0.06781245723262687

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 12.97098566666667
$ws.Cells.Item(4, 8).Value = 38.91295700000001
$ws.Cells.Item(4, 9).Value = 0.7291028508134716
$ws.Cells.Item(4, 10).Value = 0.7291028508134717
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 10.933664
$ws.Cells.Item(4, 14).Value = 32.800992
$ws.Cells.Item(4, 15).Value = 0.1034534789405002
$ws.Cells.Item(4, 16).Value = 0.1034534789405003
$ws.Cells.Item(4, 17).Value = 141.8203990281494
$ws.Cells.Item(4, 18).Value = 1276.383591253344
$ws.Cells.Item(4, 19).Value = 0.07542822642209018
$ws.Cells.Item(4, 20).Value = 0.0754282264220902

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 12.97098566666667
$ws.Cells.Item(5, 8).Value = 38.91295700000001
$ws.Cells.Item(5, 9).Value = 0.7291028508134716
$ws.Cells.Item(5, 10).Value = 0.7291028508134717
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 37.34441866666666
$ws.Cells.Item(5, 14).Value = 112.033256
$ws.Cells.Item(5, 15).Value = 0.3533499868001453
$ws.Cells.Item(5, 16).Value = 0.3533499868001453
$ws.Cells.Item(5, 17).Value = 484.3939192553325
$ws.Cells.Item(5, 18).Value = 4359.545273297992
$ws.Cells.Item(5, 19).Value = 0.2576284827108885
$ws.Cells.Item(5, 20).Value = 0.2576284827108886

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1.047813333333333
$ws.Cells.Item(6, 8).Value = 3.14344
$ws.Cells.Item(6, 9).Value = 0.05889789011308234
$ws.Cells.Item(6, 10).Value = 0.05889789011308236
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 47.57896333333333
$ws.Cells.Item(6, 14).Value = 142.73689
$ws.Cells.Item(6, 15).Value = 0.450188452948237
$ws.Cells.Item(6, 16).Value = 0.4501884529482371
$ws.Cells.Item(6, 17).Value = 49.85387216684443
$ws.Cells.Item(6, 18).Value = 448.6848495016
$ws.Cells.Item(6, 19).Value = 0.02651515003192381
$ws.Cells.Item(6, 20).Value = 0.02651515003192381

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1.047813333333333
$ws.Cells.Item(7, 8).Value = 3.14344
$ws.Cells.Item(7, 9).Value = 0.05889789011308234
$ws.Cells.Item(7, 10).Value = 0.05889789011308236
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 9.829723666666666
$ws.Cells.Item(7, 14).Value = 29.489171
$ws.Cells.Item(7, 15).Value = 0.09300808131111737
$ws.Cells.Item(7, 16).Value = 0.09300808131111739
$ws.Cells.Item(7, 17).Value = 10.29971552091556
$ws.Cells.Item(7, 18).Value = 92.69743968824
$ws.Cells.Item(7, 19).Value = 0.005477979752690819
$ws.Cells.Item(7, 20).Value = 0.005477979752690821

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1.047813333333333
$ws.Cells.Item(8, 8).Value = 3.14344
$ws.Cells.Item(8, 9).Value = 0.05889789011308234
$ws.Cells.Item(8, 10).Value = 0.05889789011308236
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 10.933664
$ws.Cells.Item(8, 14).Value = 32.800992
$ws.Cells.Item(8, 15).Value = 0.1034534789405002
$ws.Cells.Item(8, 16).Value = 0.1034534789405003
$ws.Cells.Item(8, 17).Value = 11.45643892138667
$ws.Cells.Item(8, 18).Value = 103.10795029248
$ws.Cells.Item(8, 19).Value = 0.006093191634453662
$ws.Cells.Item(8, 20).Value = 0.006093191634453664

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1.047813333333333
$ws.Cells.Item(9, 8).Value = 3.14344
$ws.Cells.Item(9, 9).Value = 0.05889789011308234
$ws.Cells.Item(9, 10).Value = 0.05889789011308236
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 37.34441866666666
$ws.Cells.Item(9, 14).Value = 112.033256
$ws.Cells.Item(9, 15).Value = 0.3533499868001453
$ws.Cells.Item(9, 16).Value = 0.3533499868001453
$ws.Cells.Item(9, 17).Value = 39.12997980451555
$ws.Cells.Item(9, 18).Value = 352.16981824064
$ws.Cells.Item(9, 19).Value = 0.02081156869401405
$ws.Cells.Item(9, 20).Value = 0.02081156869401406

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 1.312552
$ws.Cells.Item(10, 8).Value = 3.937656
$ws.Cells.Item(10, 9).Value = 0.07377892703252469
$ws.Cells.Item(10, 10).Value = 0.0737789270325247
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 47.57896333333333
$ws.Cells.Item(10, 14).Value = 142.73689
$ws.Cells.Item(10, 15).Value = 0.450188452948237
$ws.Cells.Item(10, 16).Value = 0.4501884529482371
$ws.Cells.Item(10, 17).Value = 62.44986348109332
$ws.Cells.Item(10, 18).Value = 562.04877132984
$ws.Cells.Item(10, 19).Value = 0.03321442102095315
$ws.Cells.Item(10, 20).Value = 0.03321442102095316

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 1.312552
$ws.Cells.Item(11, 8).Value = 3.937656
$ws.Cells.Item(11, 9).Value = 0.07377892703252469
$ws.Cells.Item(11, 10).Value = 0.0737789270325247
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 9.829723666666666
$ws.Cells.Item(11, 14).Value = 29.489171
$ws.Cells.Item(11, 15).Value = 0.09300808131111737
$ws.Cells.Item(11, 16).Value = 0.09300808131111739
$ws.Cells.Item(11, 17).Value = 12.90202345813067
$ws.Cells.Item(11, 18).Value = 116.118211123176
$ws.Cells.Item(11, 19).Value = 0.006862036444488052
$ws.Cells.Item(11, 20).Value = 0.006862036444488054

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 1.312552
$ws.Cells.Item(12, 8).Value = 3.937656
$ws.Cells.Item(12, 9).Value = 0.07377892703252469
$ws.Cells.Item(12, 10).Value = 0.0737789270325247
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 10.933664
$ws.Cells.Item(12, 14).Value = 32.800992
$ws.Cells.Item(12, 15).Value = 0.1034534789405002
$ws.Cells.Item(12, 16).Value = 0.1034534789405003
$ws.Cells.Item(12, 17).Value = 14.351002550528
$ws.Cells.Item(12, 18).Value = 129.159022954752
$ws.Cells.Item(12, 19).Value = 0.007632686674011997
$ws.Cells.Item(12, 20).Value = 0.007632686674011999

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 1.312552
$ws.Cells.Item(13, 8).Value = 3.937656
$ws.Cells.Item(13, 9).Value = 0.07377892703252469
$ws.Cells.Item(13, 10).Value = 0.0737789270325247
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 37.34441866666666
$ws.Cells.Item(13, 14).Value = 112.033256
$ws.Cells.Item(13, 15).Value = 0.3533499868001453
$ws.Cells.Item(13, 16).Value = 0.3533499868001453
$ws.Cells.Item(13, 17).Value = 49.01649140977066
$ws.Cells.Item(13, 18).Value = 441.148422687936
$ws.Cells.Item(13, 19).Value = 0.02606978289307148
$ws.Cells.Item(13, 20).Value = 0.02606978289307149

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 2.458986333333333
$ws.Cells.Item(14, 8).Value = 7.376958999999999
$ws.Cells.Item(14, 9).Value = 0.1382203320409214
$ws.Cells.Item(14, 10).Value = 0.1382203320409214
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 47.57896333333333
$ws.Cells.Item(14, 14).Value = 142.73689
$ws.Cells.Item(14, 15).Value = 0.450188452948237
$ws.Cells.Item(14, 16).Value = 0.4501884529482371
$ws.Cells.Item(14, 17).Value = 116.9960205908344
$ws.Cells.Item(14, 18).Value = 1052.96418531751
$ws.Cells.Item(14, 19).Value = 0.06222519744749402
$ws.Cells.Item(14, 20).Value = 0.06222519744749404

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 2.458986333333333
$ws.Cells.Item(15, 8).Value = 7.376958999999999
$ws.Cells.Item(15, 9).Value = 0.1382203320409214
$ws.Cells.Item(15, 10).Value = 0.1382203320409214
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 9.829723666666666
$ws.Cells.Item(15, 14).Value = 29.489171
$ws.Cells.Item(15, 15).Value = 0.09300808131111737
$ws.Cells.Item(15, 16).Value = 0.09300808131111739
$ws.Cells.Item(15, 17).Value = 24.17115615677655
$ws.Cells.Item(15, 18).Value = 217.540405410989
$ws.Cells.Item(15, 19).Value = 0.01285560788131166
$ws.Cells.Item(15, 20).Value = 0.01285560788131166

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 2.458986333333333
$ws.Cells.Item(16, 8).Value = 7.376958999999999
$ws.Cells.Item(16, 9).Value = 0.1382203320409214
$ws.Cells.Item(16, 10).Value = 0.1382203320409214
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 10.933664
$ws.Cells.Item(16, 14).Value = 32.800992
$ws.Cells.Item(16, 15).Value = 0.1034534789405002
$ws.Cells.Item(16, 16).Value = 0.1034534789405003
$ws.Cells.Item(16, 17).Value = 26.88573034925867
$ws.Cells.Item(16, 18).Value = 241.971573143328
$ws.Cells.Item(16, 19).Value = 0.01429937420994441
$ws.Cells.Item(16, 20).Value = 0.01429937420994441

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 2.458986333333333
$ws.Cells.Item(17, 8).Value = 7.376958999999999
$ws.Cells.Item(17, 9).Value = 0.1382203320409214
$ws.Cells.Item(17, 10).Value = 0.1382203320409214
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 37.34441866666666
$ws.Cells.Item(17, 14).Value = 112.033256
$ws.Cells.Item(17, 15).Value = 0.3533499868001453
$ws.Cells.Item(17, 16).Value = 0.3533499868001453
$ws.Cells.Item(17, 17).Value = 91.82941512761154
$ws.Cells.Item(17, 18).Value = 826.4647361485039
$ws.Cells.Item(17, 19).Value = 0.04884015250217126
$ws.Cells.Item(17, 20).Value = 0.04884015250217127
